$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.295.77'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '1.825.80'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.84%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '330.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4436'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3759'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.79'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07719'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.75%  '
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.14'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.323'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.542'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '1.832.45'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.63'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +14.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001083'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06491'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.316'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5378'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').Value = '28.361.09'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.178'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -10.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.69'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '155.58'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.349'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.68%  '
$ws.Range('D30').Value = '2.033.35'
$ws.Range('E30').Value = '  +2.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '128.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.197'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.873'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09264'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.670'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '13.04'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02344'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2182'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.179'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6579'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06186'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.210'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.129'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.03'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.392'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6083'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.775'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.048'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '126.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.153'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.28%  '
